# Fix a bug in the "local.png" figure source: slide 8's "Will push to ..."
# text box needs the remote name spelled out (.../github-collaboration/tree/...)
# and a minor wording tweak ("Then send PR on the web" -> "Then " + "send PR on the web").

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(8)
$shp = $s.Shapes.Item("文本框 38")

# Shape.Width/Height/Left/Top are expressed in points (1 pt = 12700 EMU); add
# half a point's worth of EMU head-room before dividing so the point value
# still truncates back down to the exact target EMU after the COM layer's
# internal single-precision round-trip.
function EmuToPt($emu) {
    return ($emu + 0.5) / 12700
}

$origHeightEmu = 954107
$newWidthEmu   = 5184321

# --- widen the text box to fit the extra text (cx: 4175814 -> 5184321 EMU) ---
$shp.Width = EmuToPt $newWidthEmu

$tr = $shp.TextFrame.TextRange

# --- paragraph 2: "{your new branch for PR}/tree/{PR tree name}" ---
#   -> "{your new branch for PR}" + "/" + "github" + "-collaboration" + "/tree/{PR tree name}"
$para2 = $tr.Paragraphs(2)
$existingTail = $para2.Characters(25, 21)   # "/tree/{PR tree name}"
$null = $existingTail.InsertBefore("/github-collaboration")

# Re-assert Bold on each new segment so PowerPoint keeps them as distinct runs
# (they already inherit the bold-red formatting of the run they were inserted into).
$para2.Characters(25, 1).Font.Bold = $true    # "/"
$para2.Characters(26, 6).Font.Bold = $true    # "github"
$para2.Characters(32, 14).Font.Bold = $true   # "-collaboration"

# --- paragraph 3: "Then send PR on the web" -> "Then " + "send PR on the web" ---
$para3 = $tr.Paragraphs(3)
$para3.Characters(1, 5).Font.Bold = $true     # "Then " (forces a run split)

# The text reflow above can trip this host's auto-fit simulation into growing
# the box's height; the real deck keeps the original height, so pin both
# dimensions back to their exact EMU values.
$shp.Width = EmuToPt $newWidthEmu
$shp.Height = EmuToPt $origHeightEmu
